# Generate Report for Handback
#
# - Overview sheet: "Ready for handoff" status becomes
#   "Handed back: in sync with en-US" for both rows (zh-cn + de-de columns).
# - zh-cn / de-de detail sheets: each row's "Latest Target File" (I) gets a
#   hyperlink to the source .md file, "Latest Handback File" (J) gets the
#   generated .xlf file name, and "Latest Handback DateTime" (K) gets the
#   handback timestamp.
# - Column widths widen on the affected columns to fit the new, longer text.

$wb = $excel.ActiveWorkbook

$hyperlinkBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/046daf2e279a4be68647138f045094be29932056/e2e/"
$mdName1 = "44f06573-e207-49f2-bcf6-861ad6fbb7d8.md"
$mdName2 = "713bed8e-298f-46fd-ad6d-00cc44e59c19.md"

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

# Columns E (zh-cn) and F (de-de) widen to fit the longer status text.
$wsOverview.Columns.Item(5).ColumnWidth = 29.16666666666667
$wsOverview.Columns.Item(6).ColumnWidth = 29.16666666666667

# ---------------------------------------------------------------------------
# zh-cn detail sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

$wsZh.Range("I2").Value = $mdName1
$wsZh.Range("I2").Font.Underline = 2
$wsZh.Range("I2").Font.Color = 15570276
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), ($hyperlinkBase + $mdName1), $null, $null, $mdName1)

$wsZh.Range("J2").Value = "44f06573-e207-49f2-bcf6-861ad6fbb7d8.c201ce183e39ba1be4123baf04c44233a5755220.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-10-21 01:00:06"

$wsZh.Range("I3").Value = $mdName2
$wsZh.Range("I3").Font.Underline = 2
$wsZh.Range("I3").Font.Color = 15570276
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), ($hyperlinkBase + $mdName2), $null, $null, $mdName2)

$wsZh.Range("J3").Value = "713bed8e-298f-46fd-ad6d-00cc44e59c19.01a1e8b0c5889a6e645823fa56a9d02162422ae6.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-10-21 01:00:06"

# Column C (Status), I (Latest Target File), J (Latest Handback File) widen.
$wsZh.Columns.Item(3).ColumnWidth = 29.16666666666667
$wsZh.Columns.Item(9).ColumnWidth = 39.16666666666667
$wsZh.Columns.Item(10).ColumnWidth = 39.16666666666667

# ---------------------------------------------------------------------------
# de-de detail sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$wsDe.Range("I2").Value = $mdName1
$wsDe.Range("I2").Font.Underline = 2
$wsDe.Range("I2").Font.Color = 15570276
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), ($hyperlinkBase + $mdName1), $null, $null, $mdName1)

$wsDe.Range("J2").Value = "44f06573-e207-49f2-bcf6-861ad6fbb7d8.c201ce183e39ba1be4123baf04c44233a5755220.de-de.xlf"
$wsDe.Range("K2").Value = "2016-10-21 01:00:23"

$wsDe.Range("I3").Value = $mdName2
$wsDe.Range("I3").Font.Underline = 2
$wsDe.Range("I3").Font.Color = 15570276
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), ($hyperlinkBase + $mdName2), $null, $null, $mdName2)

$wsDe.Range("J3").Value = "713bed8e-298f-46fd-ad6d-00cc44e59c19.01a1e8b0c5889a6e645823fa56a9d02162422ae6.de-de.xlf"
$wsDe.Range("K3").Value = "2016-10-21 01:00:23"

# Column C (Status), I (Latest Target File), J (Latest Handback File) widen.
$wsDe.Columns.Item(3).ColumnWidth = 29.16666666666667
$wsDe.Columns.Item(9).ColumnWidth = 39.16666666666667
$wsDe.Columns.Item(10).ColumnWidth = 39.16666666666667
